# Update cryptocurrency price/volume table to latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.864.40"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'2.285.96"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "'310.79"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").Value = "'103.14"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("D10").Value = "'38.67"
$ws.Range("E10").Value = "  -3.97%  "
$ws.Range("D11").Value = "'0.0900"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "'8.21"
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "'0.975"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'15.21"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'2.632.00"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "'2.278.67"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'42.515.31"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'7.28"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("D22").Value = "'73.15"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "'266.82"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'3.39"
$ws.Range("E24").Value = "  -4.98%  "
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "'10.75"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "'7.14"
$ws.Range("E28").Value = "  +16.92%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "'22.30"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").Value = "'35.60"
$ws.Range("E31").Value = "  -7.32%  "
$ws.Range("D32").Value = "'164.32"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "'0.0849"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "'2.56"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").Value = "'4.52"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "'0.0345"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "'2.74"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "'3.59"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").Value = "'107.00"
$ws.Range("E41").Value = "  +12.00%  "
$ws.Range("D42").Value = "'1.55"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'71.08"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "'1.715.56"
$ws.Range("E47").Value = "  +8.34%  "
$ws.Range("D48").Value = "'110.45"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "'77.11"
$ws.Range("E49").Value = "  -5.04%  "

# Rows 20/21 swap rank order (ShibaInu moves above InternetComputer)
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000104"
$ws.Range("E20").Value = "  -1.29%  "

$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'13.52"
$ws.Range("E21").Value = "  +1.83%  "

# Rows 50/51 swap rank order (THORChain moves above FraxShare)
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.15"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.62"
$ws.Range("E51").Value = "  -3.12%  "
